$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Enemigos básicos" and "Sistema de ataque (combo débil)" rows
# as completed, switching their status cell from the "pending" mark (◻️)
# to the "done" mark (✅) — the Weak Attack combo system is finished.
$ws.Range("A5").Value = "✅"
$ws.Range("A6").Value = "✅"

# Update the active cell/selection to match where the cursor ended up.
$ws.Range("B7").Select()
